$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.46387466666667
$ws.Range("H2").Value = 37.391624
$ws.Range("I2").Value = 0.635632186526332
$ws.Range("J2").Value = 0.635632186526332
$ws.Range("M2").Value = 40.25420133333333
$ws.Range("N2").Value = 120.762604
$ws.Range("O2").Value = 0.4854671023051697
$ws.Range("P2").Value = 0.4854671023051695
$ws.Range("Q2").Value = 501.7233202254328
$ws.Range("R2").Value = 4515.509882028895
$ws.Range("S2").Value = 0.3085785157248375
$ws.Range("T2").Value = 0.3085785157248374

$ws.Range("G3").Value = 12.46387466666667
$ws.Range("H3").Value = 37.391624
$ws.Range("I3").Value = 0.635632186526332
$ws.Range("J3").Value = 0.635632186526332
$ws.Range("O3").Value = 0.05577747182450057
$ws.Range("P3").Value = 0.05577747182450056
$ws.Range("Q3").Value = 57.64522091133868
$ws.Range("R3").Value = 518.806988202048
$ws.Range("S3").Value = 0.03545395637471817
$ws.Range("T3").Value = 0.03545395637471817

$ws.Range("G4").Value = 12.46387466666667
$ws.Range("H4").Value = 37.391624
$ws.Range("I4").Value = 0.635632186526332
$ws.Range("J4").Value = 0.635632186526332
$ws.Range("M4").Value = 36.68940733333334
$ws.Range("N4").Value = 110.068222
$ws.Range("O4").Value = 0.4424755596543956
$ws.Range("P4").Value = 0.4424755596543954
$ws.Range("Q4").Value = 457.2921745969477
$ws.Range("R4").Value = 4115.629571372529
$ws.Range("S4").Value = 0.2812517074675859
$ws.Range("T4").Value = 0.2812517074675858

$ws.Range("G5").Value = 12.46387466666667
$ws.Range("H5").Value = 37.391624
$ws.Range("I5").Value = 0.635632186526332
$ws.Range("J5").Value = 0.635632186526332
$ws.Range("M5").Value = 1.349902
$ws.Range("N5").Value = 4.049706
$ws.Range("O5").Value = 0.01627986621593436
$ws.Range("P5").Value = 0.01627986621593436
$ws.Range("Q5").Value = 16.82500934028267
$ws.Range("R5").Value = 151.425084062544
$ws.Range("S5").Value = 0.01034800695919052
$ws.Range("T5").Value = 0.01034800695919052

$ws.Range("I6").Value = 0.1376366783586857
$ws.Range("J6").Value = 0.1376366783586857
$ws.Range("M6").Value = 40.25420133333333
$ws.Range("N6").Value = 120.762604
$ws.Range("O6").Value = 0.4854671023051697
$ws.Range("P6").Value = 0.4854671023051695
$ws.Range("Q6").Value = 108.6407087537551
$ws.Range("R6").Value = 977.766378783796
$ws.Range("S6").Value = 0.06681807941369981
$ws.Range("T6").Value = 0.0668180794136998

$ws.Range("I7").Value = 0.1376366783586857
$ws.Range("J7").Value = 0.1376366783586857
$ws.Range("O7").Value = 0.05577747182450057
$ws.Range("P7").Value = 0.05577747182450056
$ws.Range("S7").Value = 0.007677025949169441
$ws.Range("T7").Value = 0.00767702594916944

$ws.Range("I8").Value = 0.1376366783586857
$ws.Range("J8").Value = 0.1376366783586857
$ws.Range("M8").Value = 36.68940733333334
$ws.Range("N8").Value = 110.068222
$ws.Range("O8").Value = 0.4424755596543956
$ws.Range("P8").Value = 0.4424755596543954
$ws.Range("Q8").Value = 99.01980624188649
$ws.Range("R8").Value = 891.1782561769783
$ws.Range("S8").Value = 0.06090086628573151
$ws.Range("T8").Value = 0.06090086628573149

$ws.Range("I9").Value = 0.1376366783586857
$ws.Range("J9").Value = 0.1376366783586857
$ws.Range("M9").Value = 1.349902
$ws.Range("N9").Value = 4.049706
$ws.Range("O9").Value = 0.01627986621593436
$ws.Range("P9").Value = 0.01627986621593436
$ws.Range("Q9").Value = 3.643205061099334
$ws.Range("R9").Value = 32.78884554989401
$ws.Range("S9").Value = 0.002240706710084992
$ws.Range("T9").Value = 0.002240706710084992

$ws.Range("G10").Value = 3.427721
$ws.Range("H10").Value = 10.283163
$ws.Range("I10").Value = 0.1748067797776496
$ws.Range("J10").Value = 0.1748067797776496
$ws.Range("M10").Value = 40.25420133333333
$ws.Range("N10").Value = 120.762604
$ws.Range("O10").Value = 0.4854671023051697
$ws.Range("P10").Value = 0.4854671023051695
$ws.Range("Q10").Value = 137.9801712484947
$ws.Range("R10").Value = 1241.821541236452
$ws.Range("S10").Value = 0.08486294084195345
$ws.Range("T10").Value = 0.08486294084195344

$ws.Range("G11").Value = 3.427721
$ws.Range("H11").Value = 10.283163
$ws.Range("I11").Value = 0.1748067797776496
$ws.Range("J11").Value = 0.1748067797776496
$ws.Range("O11").Value = 0.05577747182450057
$ws.Range("P11").Value = 0.05577747182450056
$ws.Range("Q11").Value = 15.853154781464
$ws.Range("R11").Value = 142.678393033176
$ws.Range("S11").Value = 0.009750280233779524
$ws.Range("T11").Value = 0.009750280233779522

$ws.Range("G12").Value = 3.427721
$ws.Range("H12").Value = 10.283163
$ws.Range("I12").Value = 0.1748067797776496
$ws.Range("J12").Value = 0.1748067797776496
$ws.Range("M12").Value = 36.68940733333334
$ws.Range("N12").Value = 110.068222
$ws.Range("O12").Value = 0.4424755596543956
$ws.Range("P12").Value = 0.4424755596543954
$ws.Range("Q12").Value = 125.7610519940207
$ws.Range("R12").Value = 1131.849467946186
$ws.Range("S12").Value = 0.07734772771349815
$ws.Range("T12").Value = 0.07734772771349814

$ws.Range("G13").Value = 3.427721
$ws.Range("H13").Value = 10.283163
$ws.Range("I13").Value = 0.1748067797776496
$ws.Range("J13").Value = 0.1748067797776496
$ws.Range("M13").Value = 1.349902
$ws.Range("N13").Value = 4.049706
$ws.Range("O13").Value = 0.01627986621593436
$ws.Range("P13").Value = 0.01627986621593436
$ws.Range("Q13").Value = 4.627087433342001
$ws.Range("R13").Value = 41.64378690007801
$ws.Range("S13").Value = 0.002845830988418435
$ws.Range("T13").Value = 0.002845830988418434

$ws.Range("G14").Value = 1.018165333333333
$ws.Range("H14").Value = 3.054496
$ws.Range("I14").Value = 0.05192435533733264
$ws.Range("J14").Value = 0.05192435533733263
$ws.Range("M14").Value = 40.25420133333333
$ws.Range("N14").Value = 120.762604
$ws.Range("O14").Value = 0.4854671023051697
$ws.Range("P14").Value = 0.4854671023051695
$ws.Range("Q14").Value = 40.98543231862044
$ws.Range("R14").Value = 368.868890867584
$ws.Range("S14").Value = 0.02520756632467885
$ws.Range("T14").Value = 0.02520756632467884

$ws.Range("G15").Value = 1.018165333333333
$ws.Range("H15").Value = 3.054496
$ws.Range("I15").Value = 0.05192435533733264
$ws.Range("J15").Value = 0.05192435533733263
$ws.Range("O15").Value = 0.05577747182450057
$ws.Range("P15").Value = 0.05577747182450056
$ws.Range("Q15").Value = 4.708998376021334
$ws.Range("R15").Value = 42.380985384192
$ws.Range("S15").Value = 0.002896209266833427
$ws.Range("T15").Value = 0.002896209266833426

$ws.Range("G16").Value = 1.018165333333333
$ws.Range("H16").Value = 3.054496
$ws.Range("I16").Value = 0.05192435533733264
$ws.Range("J16").Value = 0.05192435533733263
$ws.Range("M16").Value = 36.68940733333334
$ws.Range("N16").Value = 110.068222
$ws.Range("O16").Value = 0.4424755596543956
$ws.Range("P16").Value = 0.4424755596543954
$ws.Range("Q16").Value = 37.35588264734579
$ws.Range("R16").Value = 336.2029438261121
$ws.Range("S16").Value = 0.02297525818757996
$ws.Range("T16").Value = 0.02297525818757995

$ws.Range("G17").Value = 1.018165333333333
$ws.Range("H17").Value = 3.054496
$ws.Range("I17").Value = 0.05192435533733264
$ws.Range("J17").Value = 0.05192435533733263
$ws.Range("M17").Value = 1.349902
$ws.Range("N17").Value = 4.049706
$ws.Range("O17").Value = 0.01627986621593436
$ws.Range("P17").Value = 0.01627986621593436
$ws.Range("Q17").Value = 1.374423419797334
$ws.Range("R17").Value = 12.369810778176
$ws.Range("S17").Value = 0.0008453215582404127
$ws.Range("T17").Value = 0.0008453215582404123
